$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "29.221.13"
$ws.Range("E2").Value2 = "  +0.11%  "

$ws.Range("D3").Value2 = "1.857.11"
$ws.Range("E3").Value2 = "  +0.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value2 = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "0.7105"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = "  +0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "237.71"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = "  -0.54%  "

$ws.Range("E7").Value2 = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.08162"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value2 = "  +9.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.3039"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value2 = "  -0.48%  "

$ws.Range("E10").Value2 = "  -0.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.08200"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value2 = "  +0.59%  "

$ws.Range("D12").Value2 = "1.870.44"
$ws.Range("E12").Value2 = "  +0.46%  "

$ws.Range("E13").Value2 = "  -0.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "0.7078"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value2 = "  -2.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "89.50"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value2 = "  +0.36%  "

$ws.Range("D16").Value2 = "29.267.12"
$ws.Range("E16").Value2 = "  -0.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "0.000007913"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value2 = "  +3.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "5.790"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value2 = "  +0.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "13.34"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value2 = "  +1.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "237.83"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value2 = "  -0.48%  "

$ws.Range("B21").Value2 = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value2 = "2.155.78"
$ws.Range("E21").Value2 = "  +0.71%  "

$ws.Range("B22").Value2 = "Dai"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "1.001"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value2 = "  +0.03%  "

$ws.Range("E23").Value2 = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "7.420"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value2 = "  -2.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "162.54"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value2 = "  +1.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "0.1463"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value2 = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "8.953"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value2 = "  -0.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "18.08"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value2 = "  -0.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "1.956"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value2 = "  -1.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.427"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value2 = "  +1.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "4.401"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value2 = "  -2.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "1.480"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value2 = "  -1.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "4.022"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value2 = "  +0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "0.05221"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value2 = "  +0.27%  "

$ws.Range("E35").Value2 = "  -1.76%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.7081"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value2 = "  +0.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.9990"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value2 = "  -3.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "2.673"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value2 = "  +0.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.01857"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value2 = "  -0.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "2.728"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value2 = "  +1.83%  "

$ws.Range("B41").Value2 = "TrustWalletToken"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "0.9228"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value2 = "  -2.16%  "

$ws.Range("B42").Value2 = "Maker"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value2 = "1.139.36"
$ws.Range("E42").Value2 = "  +6.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.4283"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value2 = "  -0.63%  "

$ws.Range("E44").Value2 = "  -2.41%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "70.11"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value2 = "  -0.55%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "0.9992"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value2 = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "102.34"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = "  -1.23%  "

$ws.Range("E48").Value2 = "  +1.57%  "

$ws.Range("D49").Value2 = "2.010.95"
$ws.Range("E49").Value2 = "  -1.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "9.189"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value2 = "  +0.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "6.969"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value2 = "  -1.20%  "
